$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values are stored as text (matches source inlineStr "t" type),
# by marking the cells as Text-formatted before assigning the string value.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '244.86'
$ws.Range("D3").Value = '21.82'
$ws.Range("D4").Value = '5.384'
$ws.Range("D5").Value = '0.06021'
$ws.Range("D7").Value = '0.8153'
$ws.Range("D8").Value = '0.9320'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '0.1433'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '0.07415'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.03456'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03068'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09407'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = '4.005'
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001601'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '0.04796'
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '0.0005941'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("D18").Value = '0.005569'
$ws.Range("D19").Value = '0.004156'
$ws.Range("D20").Value = '0.0009905'
$ws.Range("D21").Value = '3.668'
$ws.Range("D22").Value = '6.425'
$ws.Range("D23").Value = '2.172'
$ws.Range("D25").Value = '0.1339'
$ws.Range("D26").Value = '0.00007000'
$ws.Range("D40").Value = '0.04016'
$ws.Range("D41").Value = '0.006420'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("D42").Value = '0.1075'
$ws.Range("D43").Value = '0.002720'
$ws.Range("D44").Value = '0.006659'
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("D45").Value = '0.00005270'
$ws.Range("D47").Value = '0.8602'
$ws.Range("D48").Value = '0.002520'
